# Adds units (uot) to the time-series header strings in row 1, and
# refreshes the column widths / active-cell selection to match how the
# sheet looked after the author widened the header columns by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("import")

# --- Append unit suffixes to the two series-name headers (B1, C1) ---
$ws.Range("B1").Value = "02600.Flow.Inst.~1Day.0.DailyComputed[CFS] "
$ws.Range("C1").Value = "01080.Stage.Inst.1Day.0.Manual 0700[FEET]"

# --- Widen B & C so the longer header text fits, and drop the old
#     "best fit" auto-sizing in favour of an explicit custom width ---
$ws.Columns.Item(2).ColumnWidth = 37.648995535714285
$ws.Columns.Item(3).ColumnWidth = 41.738839285714285

# --- Move the active selection from B3 to C1 ---
$null = $ws.Range("C1").Select()
